# Update metrics on "Dados dos testes" sheet (row 2) with new values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados dos testes")

$ws.Range("B2").Value = 0.9212
$ws.Range("C2").Value = 0.9301
$ws.Range("D2").Value = 0.9187
$ws.Range("E2").Value = 0.9196
$ws.Range("F2").Value = 0.9407
$ws.Range("G2").Value = 0.8966
$ws.Range("H2").Value = 0.1034
$ws.Range("I2").Value = 0.0593
$ws.Range("J2").Value = 373
$ws.Range("K2").Value = 43
$ws.Range("L2").Value = 31
$ws.Range("M2").Value = 492
$ws.Range("N2").Value = 939
